# Revert config file handling
# Append a new row (45) to each of the four data sheets, duplicating the
# last existing row (44) with an updated timestamp, matching the pattern
# of the daily-snapshot rows already in each sheet. A few sheets also get
# small value tweaks on top of the duplicated row.

$wb = $excel.ActiveWorkbook

$newTimestamp = 45831.4970949074
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- Sheet "FE_LFT_#1": plain duplicate of row 44 with new timestamp ---
$ws = $wb.Worksheets.Item("FE_LFT_#1")
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial()
$ws.Range("A45").Value = $newTimestamp
$ws.Range("A45").NumberFormat = $dateFormat

# --- Sheet "FE_LFT_#2": duplicate row 44, then adjust D45/H45 ---
$ws = $wb.Worksheets.Item("FE_LFT_#2")
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial()
$ws.Range("A45").Value = $newTimestamp
$ws.Range("A45").NumberFormat = $dateFormat
$ws.Range("D45").Value = "0x01,0x70"
$ws.Range("H45").Value = 368

# --- Sheet "FE_PLT_#1": plain duplicate of row 44 with new timestamp ---
$ws = $wb.Worksheets.Item("FE_PLT_#1")
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial()
$ws.Range("A45").Value = $newTimestamp
$ws.Range("A45").NumberFormat = $dateFormat

# --- Sheet "FE_PLT_#2": duplicate row 44, then adjust D45/H45 ---
$ws = $wb.Worksheets.Item("FE_PLT_#2")
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial()
$ws.Range("A45").Value = $newTimestamp
$ws.Range("A45").NumberFormat = $dateFormat
$ws.Range("D45").Value = "0x00,0x69"
$ws.Range("H45").Value = 105
